$d = $word.ActiveDocument

function Find-ParagraphIndex($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$needle*") {
            return $i
        }
    }
    return -1
}

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# Locate all three target paragraphs BEFORE any edits are made, since the
# first edit introduces a new literal "{no2_1_5_2}" substring elsewhere in
# the document which would otherwise make later text searches ambiguous.
# Range objects captured now stay correctly anchored to their paragraph even
# after earlier parts of the document are edited.
$idx1 = Find-ParagraphIndex("no2_1_5_1")
$idx2 = Find-ParagraphIndex("no2_1_5_2")
$idx3 = Find-ParagraphIndex("ให้บริการครอบคลุม")

$r1 = $d.Paragraphs.Item($idx1).Range
$r2 = $d.Paragraphs.Item($idx2).Range
$r3 = $d.Paragraphs.Item($idx3).Range

# ---------------------------------------------------------------------------
# Change 1: paragraph that holds "...{no2_1_5_1}" gains the paragraph-mark
# w:hint="cs"/w:cs toggle, and 4 new runs (space, "จำนวน", " {no2_1_5_2} ",
# "สาย") get appended after the {no2_1_5_1} run.
# ---------------------------------------------------------------------------
$body1 = '<w:body><w:p>' +
    '<w:pPr><w:rPr><w:rFonts w:ascii="TH SarabunPSK" w:eastAsia="Times New Roman" w:hAnsi="TH SarabunPSK" w:cs="TH SarabunPSK" w:hint="cs"/><w:cs/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="TH SarabunPSK" w:eastAsia="Times New Roman" w:hAnsi="TH SarabunPSK" w:cs="TH SarabunPSK"/><w:cs/></w:rPr><w:t>อื่นๆ (ระบุ)</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="TH SarabunPSK" w:eastAsia="Times New Roman" w:hAnsi="TH SarabunPSK" w:cs="TH SarabunPSK"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="TH SarabunPSK" w:eastAsia="Times New Roman" w:hAnsi="TH SarabunPSK" w:cs="TH SarabunPSK"/></w:rPr><w:t>{no2_1_5_1}</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="TH SarabunPSK" w:eastAsia="Times New Roman" w:hAnsi="TH SarabunPSK" w:cs="TH SarabunPSK"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="TH SarabunPSK" w:eastAsia="Times New Roman" w:hAnsi="TH SarabunPSK" w:cs="TH SarabunPSK" w:hint="cs"/><w:cs/></w:rPr><w:t>จำนวน</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="TH SarabunPSK" w:eastAsia="Times New Roman" w:hAnsi="TH SarabunPSK" w:cs="TH SarabunPSK"/></w:rPr><w:t xml:space="preserve"> {no2_1_5_2} </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="TH SarabunPSK" w:eastAsia="Times New Roman" w:hAnsi="TH SarabunPSK" w:cs="TH SarabunPSK" w:hint="cs"/><w:cs/></w:rPr><w:t>สาย</w:t></w:r>' +
    '</w:p></w:body>'
$r1.InsertXML($pkgHeader + $body1 + $pkgFooter)

# ---------------------------------------------------------------------------
# Change 2: paragraph that holds "{no2_1_5_2}" (the field placeholder cell)
# changes text to "{no2_1_6}" and the _GoBack bookmark moves to just after
# that run.
# ---------------------------------------------------------------------------
$body2 = '<w:body><w:p>' +
    '<w:pPr><w:rPr><w:rFonts w:ascii="TH SarabunPSK" w:eastAsia="Times New Roman" w:hAnsi="TH SarabunPSK" w:cs="TH SarabunPSK"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="TH SarabunPSK" w:eastAsia="Times New Roman" w:hAnsi="TH SarabunPSK" w:cs="TH SarabunPSK"/></w:rPr><w:t>{no2_1_6}</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '</w:p></w:body>'
$r2.InsertXML($pkgHeader + $body2 + $pkgFooter)

# ---------------------------------------------------------------------------
# Change 3: paragraph "ให้บริการครอบคลุม " + bookmark(_GoBack) + "ร้อยละ" + " "
# becomes a merged run "ให้บริการครอบคลุม ร้อยละ" + " ", bookmark removed
# (it moved to change 2's location).
# ---------------------------------------------------------------------------
$body3 = '<w:body><w:p>' +
    '<w:pPr><w:rPr><w:rFonts w:ascii="TH SarabunPSK" w:eastAsia="Times New Roman" w:hAnsi="TH SarabunPSK" w:cs="TH SarabunPSK"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="TH SarabunPSK" w:eastAsia="Times New Roman" w:hAnsi="TH SarabunPSK" w:cs="TH SarabunPSK" w:hint="cs"/><w:cs/></w:rPr><w:t>ให้บริการครอบคลุม ร้อยละ</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="TH SarabunPSK" w:eastAsia="Times New Roman" w:hAnsi="TH SarabunPSK" w:cs="TH SarabunPSK"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '</w:p></w:body>'
$r3.InsertXML($pkgHeader + $body3 + $pkgFooter)

Write-Host "Done"
